# Update the cryptocurrency symbol list with freshly scraped values.
# (Source values are stored as text in the sheet, e.g. "236.36", not as
# numbers, so we prefix with an apostrophe to force text entry and then
# clear the resulting "quote-prefixed" formatting so the cell style stays
# the same as the original - matching how the data was originally written.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "236.36"
    "D3"  = "21.72"
    "D4"  = "5.356"
    "D5"  = "0.05588"
    "D6"  = "6.467"
    "D7"  = "3.356"
    "D8"  = "0.7992"
    "D9"  = "1.036"
    "D10" = "0.1391"
    "D11" = "0.07324"
    "D12" = "0.03162"
    "D13" = "0.02977"
    "D14" = "0.09248"
    "D16" = "3.255"
    "D17" = "0.04791"
    "E18" = "17OneONE"
    "D19" = "0.006222"
    "D20" = "0.005041"
    "D21" = "0.001053"
    "D22" = "0.0001502"
    "D23" = "0.0003998"
    "D24" = "3.951"
    "D27" = "0.1044"
    "D40" = "0.04112"
    "D41" = "0.007023"
    "D42" = "0.003505"
    "D43" = "0.1033"
    "D44" = "0.008807"
    "D45" = "0.00005444"
    "D47" = "0.6760"
    "D48" = "0.03487"
    "E48" = "47BOLOBOLOWorstin24h"
    "D49" = "0.00002103"
    "D50" = "0.01011"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.ClearFormats()
}
